# Issue #11 Exception on directory images
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Issues")

# Row 12: new issue #11
$ws.Cells.Item(12, 1).Value = 11
$ws.Cells.Item(12, 2).Value = "DONE"
$ws.Cells.Item(12, 3).Value = "Server"
$ws.Cells.Item(12, 6).Value = "Exception on playing images in directories"
$ws.Cells.Item(12, 4).Value = "Exception on directory images"
$ws.Range("A12:F12").RowHeight = 29

# Row 13: new issue #12
$ws.Cells.Item(13, 1).Value = 12
$ws.Cells.Item(13, 3).Value = "UI"
$ws.Cells.Item(13, 4).Value = "Better navigation from edit to play"
$ws.Cells.Item(13, 6).Value = "Better navigation from edit to play"
$ws.Range("A13:F13").RowHeight = 29

$ws.Range("B12").Select()
